# Reln-Itgb1 NATMI output refreshed with new TPM inputs:
#  - Rows 2-13 (ECs/FAPs/MuSCs senders) keep the same sending/ligand/receptor/
#    target cluster labels, but every computed metric (detection rates,
#    expression values, specificities, edge weights) is recalculated.
#  - The "Resolving-Mac" sending-cluster block (old rows 14-17) is dropped
#    entirely now that it no longer has ligand-expressing cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recomputed numeric values for rows 2-13 (columns E:T) per new TPM inputs
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1471086666666667
$ws.Range("H2").Value = 0.441326
$ws.Range("I2").Value = 0.03503939655440032
$ws.Range("J2").Value = 0.03503939655440032
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 21.43384228025378
$ws.Range("R2").Value = 192.904580522284
$ws.Range("S2").Value = 0.01004207293261845
$ws.Range("T2").Value = 0.01004207293261845
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1471086666666667
$ws.Range("H3").Value = 0.441326
$ws.Range("I3").Value = 0.03503939655440032
$ws.Range("J3").Value = 0.03503939655440032
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 24.83189919302312
$ws.Range("R3").Value = 223.487092737208
$ws.Range("S3").Value = 0.01163411298316294
$ws.Range("T3").Value = 0.01163411298316294
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1471086666666667
$ws.Range("H4").Value = 0.441326
$ws.Range("I4").Value = 0.03503939655440032
$ws.Range("J4").Value = 0.03503939655440032
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 18.84846777814045
$ws.Range("R4").Value = 169.636210003264
$ws.Range("S4").Value = 0.008830786642046432
$ws.Range("T4").Value = 0.008830786642046432
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1471086666666667
$ws.Range("H5").Value = 0.441326
$ws.Range("I5").Value = 0.03503939655440032
$ws.Range("J5").Value = 0.03503939655440032
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 9.674024650251333
$ws.Range("R5").Value = 87.06622185226199
$ws.Range("S5").Value = 0.004532423996572503
$ws.Range("T5").Value = 0.004532423996572503
$ws.Range("I6").Value = 0.2452691860358486
$ws.Range("J6").Value = 0.2452691860358485
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 150.0328649078407
$ws.Range("R6").Value = 1350.295784170566
$ws.Range("S6").Value = 0.07029262191978712
$ws.Range("T6").Value = 0.07029262191978711
$ws.Range("I7").Value = 0.2452691860358486
$ws.Range("J7").Value = 0.2452691860358485
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("S7").Value = 0.08143660285927858
$ws.Range("T7").Value = 0.08143660285927855
$ws.Range("I8").Value = 0.2452691860358486
$ws.Range("J8").Value = 0.2452691860358485
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 131.9357296233707
$ws.Range("R8").Value = 1187.421566610336
$ws.Range("S8").Value = 0.06181384569189941
$ws.Range("T8").Value = 0.0618138456918994
$ws.Range("I9").Value = 0.2452691860358486
$ws.Range("J9").Value = 0.2452691860358485
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 67.71635316190699
$ws.Range("R9").Value = 609.4471784571629
$ws.Range("S9").Value = 0.03172611556488351
$ws.Range("T9").Value = 0.03172611556488351
$ws.Range("G10").Value = 3.021537333333333
$ws.Range("H10").Value = 9.064612
$ws.Range("I10").Value = 0.7196914174097511
$ws.Range("J10").Value = 0.7196914174097511
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 440.2402395048009
$ws.Range("R10").Value = 3962.162155543209
$ws.Range("S10").Value = 0.2062590801581787
$ws.Range("T10").Value = 0.2062590801581787
$ws.Range("G11").Value = 3.021537333333333
$ws.Range("H11").Value = 9.064612
$ws.Range("I11").Value = 0.7196914174097511
$ws.Range("J11").Value = 0.7196914174097511
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 510.0346034628996
$ws.Range("R11").Value = 4590.311431166097
$ws.Range("S11").Value = 0.2389587745941426
$ws.Range("T11").Value = 0.2389587745941426
$ws.Range("G12").Value = 3.021537333333333
$ws.Range("H12").Value = 9.064612
$ws.Range("I12").Value = 0.7196914174097511
$ws.Range("J12").Value = 0.7196914174097511
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 387.1379597017743
$ws.Range("R12").Value = 3484.241637315968
$ws.Range("S12").Value = 0.1813798746616646
$ws.Range("T12").Value = 0.1813798746616646
$ws.Range("G13").Value = 3.021537333333333
$ws.Range("H13").Value = 9.064612
$ws.Range("I13").Value = 0.7196914174097511
$ws.Range("J13").Value = 0.7196914174097511
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 198.6995552787826
$ws.Range("R13").Value = 1788.295997509044
$ws.Range("S13").Value = 0.0930936879957652
$ws.Range("T13").Value = 0.0930936879957652

# Remove the obsolete "Resolving-Mac" sending-cluster block (rows 14-17)
$ws.Rows("14:17").Delete()
